# Update LR-pair TPM metrics on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "I2" = 0.9223046214701632
    "J2" = 0.9223046214701632
    "K2" = 3
    "L2" = 1
    "M2" = 1.315861666666667
    "N2" = 3.947585
    "O2" = 0.2754050739440597
    "P2" = 0.2754050739440597
    "Q2" = 2.195542823928334
    "R2" = 19.759885415355
    "S2" = 0.2540073724749383
    "T2" = 0.2540073724749383

    "I3" = 0.9223046214701632
    "J3" = 0.9223046214701632
    "O3" = 0.3040809095127364
    "P3" = 0.3040809095127364
    "S3" = 0.2804552281444473
    "T3" = 0.2804552281444473

    "I4" = 0.9223046214701632
    "J4" = 0.9223046214701632
    "M4" = 2.009179666666667
    "N4" = 6.027539
    "O4" = 0.4205140165432039
    "P4" = 0.4205140165432039
    "Q4" = 3.352358466606334
    "R4" = 30.171226199457
    "S4" = 0.3878420208507776
    "T4" = 0.3878420208507776

    "E5" = 3
    "F5" = 1
    "G5" = 0.140557
    "H5" = 0.421671
    "I5" = 0.07769537852983674
    "J5" = 0.07769537852983674
    "K5" = 3
    "L5" = 1
    "M5" = 1.315861666666667
    "N5" = 3.947585
    "O5" = 0.2754050739440597
    "P5" = 0.2754050739440597
    "Q5" = 0.1849535682816667
    "R5" = 1.664582114535
    "S5" = 0.0213977014691214
    "T5" = 0.0213977014691214

    "E6" = 3
    "F6" = 1
    "G6" = 0.140557
    "H6" = 0.421671
    "I6" = 0.07769537852983674
    "J6" = 0.07769537852983674
    "O6" = 0.3040809095127364
    "P6" = 0.3040809095127364
    "Q6" = 0.2042113765563333
    "R6" = 1.837902389007
    "S6" = 0.02362568136828909
    "T6" = 0.02362568136828909

    "E7" = 3
    "F7" = 1
    "G7" = 0.140557
    "H7" = 0.421671
    "I7" = 0.07769537852983674
    "J7" = 0.07769537852983674
    "M7" = 2.009179666666667
    "N7" = 6.027539
    "O7" = 0.4205140165432039
    "P7" = 0.4205140165432039
    "Q7" = 0.2824042664076667
    "R7" = 2.541638397669
    "S7" = 0.03267199569242626
    "T7" = 0.03267199569242626
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
